$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.121.15'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '1.919.51'
$ws.Range("E3").Value = '  +2.57%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "'319.36"
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "'0.5074"
$ws.Range("E7").Value = '  -0.40%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "'0.4064"
$ws.Range("E8").Value = '  +2.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "'0.08336"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "'1.116"
$ws.Range("E10").Value = '  +2.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "'42.08"
$ws.Range("E11").Value = '  -0.25%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "'24.21"
$ws.Range("E12").Value = '  +2.81%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "'6.417"
$ws.Range("E13").Value = '  +1.85%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.919.46'
$ws.Range("E14").Value = '  +2.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "'7.256"
$ws.Range("E15").Value = '  +0.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "'1.002"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "'92.61"
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "'0.00001095"
$ws.Range("E18").Value = '  +0.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "'0.06498"
$ws.Range("E19").Value = '  +1.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "'18.50"
$ws.Range("E20").Value = '  +3.06%  '

$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "'5.947"
$ws.Range("E22").Value = '  +1.88%  '

$ws.Range("D23").Value = '30.121.60'
$ws.Range("E23").Value = '  +0.47%  '

$ws.Range("E24").Value = '  +2.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "'2.198"
$ws.Range("E25").Value = '  +1.09%  '

$ws.Range("D26").Value = '2.137.26'
$ws.Range("E26").Value = '  +2.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "'21.84"
$ws.Range("E27").Value = '  +3.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "'162.46"
$ws.Range("E28").Value = '  +1.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "'2.265"
$ws.Range("E29").Value = '  +1.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "'128.98"
$ws.Range("E30").Value = '  +0.99%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "'1.136"
$ws.Range("E31").Value = '  +6.26%  '

$ws.Range("E32").Value = '  +1.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "'5.945"
$ws.Range("E33").Value = '  +0.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "'3.795"
$ws.Range("E34").Value = '  +2.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "'0.02447"
$ws.Range("E35").Value = '  +0.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "'5.316"
$ws.Range("E36").Value = '  +1.90%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "'0.06428"
$ws.Range("E37").Value = '  +0.96%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "'1.223"
$ws.Range("E38").Value = '  +3.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "'0.2146"
$ws.Range("E39").Value = '  +0.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "'0.6466"
$ws.Range("E40").Value = '  +2.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "'8.596"
$ws.Range("E41").Value = '  +0.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "'11.47"
$ws.Range("E42").Value = '  +0.97%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "'1.213"
$ws.Range("E43").Value = '  +0.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "'13.31"
$ws.Range("E44").Value = '  +3.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "'0.6048"
$ws.Range("E45").Value = '  +2.33%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "'2.175"
$ws.Range("E46").Value = '  +8.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "'3.624"
$ws.Range("E47").Value = '  -0.46%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "'122.28"
$ws.Range("E48").Value = '  -0.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "'1.208"

$ws.Range("E50").Value = '  +1.06%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "'78.08"
$ws.Range("E51").Value = '  +1.14%  '
